$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.137.94"
$ws.Range("E2").Value = "  +0.25%  "

$ws.Range("D3").Value = "2.326.88"

$ws.Range("E4").Value = "  +0.16%  "

$ws.Range("D5").Value = "528.49"
$ws.Range("E5").Value = "  +1.92%  "

$ws.Range("D6").Value = "132.16"
$ws.Range("E6").Value = "  -2.56%  "

$ws.Range("D8").Value = "0.534"
$ws.Range("E8").Value = "  -0.58%  "

$ws.Range("D9").Value = "2.355.29"
$ws.Range("E9").Value = "  +0.12%  "

$ws.Range("E10").Value = "  -1.33%  "

$ws.Range("E11").Value = "  +0.47%  "

$ws.Range("D12").Value = "5.30"
$ws.Range("E12").Value = "  -1.96%  "

$ws.Range("D13").Value = "0.344"
$ws.Range("E13").Value = "  +0.18%  "

$ws.Range("D14").Value = "2.747.21"
$ws.Range("E14").Value = "  -0.50%  "

$ws.Range("D15").Value = "23.53"
$ws.Range("E15").Value = "  -1.84%  "

$ws.Range("D16").Value = "57.169.25"
$ws.Range("E16").Value = "  +0.41%  "

$ws.Range("D18").Value = "2.341.56"
$ws.Range("E18").Value = "  -0.71%  "

$ws.Range("D19").Value = "336.75"
$ws.Range("E19").Value = "  +2.99%  "

$ws.Range("D20").Value = "10.46"
$ws.Range("E20").Value = "  -1.30%  "

$ws.Range("D21").Value = "6.94"
$ws.Range("E21").Value = "  +2.99%  "

$ws.Range("D22").Value = "4.17"
$ws.Range("E22").Value = "  -1.51%  "

$ws.Range("E23").Value = "  +0.04%  "

$ws.Range("D24").Value = "61.76"
$ws.Range("E24").Value = "  +1.16%  "

$ws.Range("D25").Value = "8.88"
$ws.Range("E25").Value = "  +11.24%  "

$ws.Range("D26").Value = "0.165"
$ws.Range("E26").Value = "  -0.02%  "

$ws.Range("D27").Value = "0.990"
$ws.Range("E27").Value = "  -0.65%  "

$ws.Range("E28").Value = "  +3.07%  "

$ws.Range("D29").Value = "169.84"
$ws.Range("E29").Value = "  -0.49%  "

$ws.Range("D30").Value = "1.72"
$ws.Range("E30").Value = "  +0.98%  "

$ws.Range("D31").Value = "0.0₃0725"
$ws.Range("E31").Value = "  -2.42%  "

$ws.Range("E32").Value = "  -1.62%  "

$ws.Range("D33").Value = "18.55"
$ws.Range("E33").Value = "  -0.01%  "

$ws.Range("E35").Value = "  -0.31%  "

$ws.Range("D36").Value = "1.27"
$ws.Range("E36").Value = "  -0.51%  "

$ws.Range("D37").Value = "4.03"
$ws.Range("E37").Value = "  +0.05%  "

$ws.Range("D38").Value = "0.912"
$ws.Range("E38").Value = "  -0.38%  "

$ws.Range("E39").Value = "  +1.11%  "

$ws.Range("D40").Value = "38.96"
$ws.Range("E40").Value = "  +1.33%  "

$ws.Range("D41").Value = "148.31"
$ws.Range("E41").Value = "  +0.10%  "

$ws.Range("E42").Value = "  -1.08%  "

$ws.Range("B43").Value = "Filecoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D43").Value = "3.60"
$ws.Range("E43").Value = "  -1.20%  "

$ws.Range("B44").Value = "Bittensor"
$ws.Range("C44").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D44").Value = "284.82"
$ws.Range("E44").Value = "  +1.64%  "

$ws.Range("D45").Value = "5.13"
$ws.Range("E45").Value = "  -2.76%  "

$ws.Range("E46").Value = "  -0.30%  "

$ws.Range("E47").Value = "  -0.57%  "

$ws.Range("E48").Value = "  -0.27%  "

$ws.Range("D49").Value = "18.73"
$ws.Range("E49").Value = "  +4.27%  "

$ws.Range("E50").Value = "  -0.93%  "

$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "17.38"
$ws.Range("E51").Value = "  -0.35%  "
